$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B41 value from 967 to 699
$ws.Range("B41").Value = 699

# Add row 42 values: B42 = 718, C42 = 85
$ws.Range("B42").Value = 718
$ws.Range("C42").Value = 85

# Update the selection (active cell) to C43
$ws.Range("C43").Select()
